# DOMA-3100 add formatter convert to number for some colomns
#
# The ticket-analytics export template has two data rows (row 2 = "{d.tickets[i]...}",
# row 3 = "{d.tickets[i+1]...}"). Columns C:H on those rows hold the numeric ticket
# counters (processing / completed / canceled / deferred / closed / new_or_reopened).
# This adds the docxtemplater ":formatN()" number formatter to each of those
# placeholders and switches the cells' display number format to an integer ("0")
# so the exported numbers render as numbers instead of plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C", "D", "E", "F", "G", "H")
$rows = @(2, 3)

foreach ($r in $rows) {
    foreach ($col in $cols) {
        $addr = "$col$r"
        $cell = $ws.Range($addr)
        $current = $cell.Value2
        $cell.Value2 = ($current -replace '\}\s*$', ':formatN()}')
        $cell.NumberFormat = "0"
    }
}
